# Updating indexer lib data:
# Fill in the remaining 2023 months (Aug-Dec) on row 2, and move the
# "Fonte" source-link cell from O3 up to O2 (removing it from O3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (year 2023): add Agosto(H) .. Dezembro(L) monthly values.
$ws.Range("H2").Value = 1.07
$ws.Range("I2").Value = 1.1399999999999999
$ws.Range("J2").Value = 0.97
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.92

# Move the source URL from O3 to O2, matching O3's original (unstyled) look.
$ws.Range("O2").ClearFormats()
$ws.Range("O2").Value = $ws.Range("O3").Value2
$ws.Range("O3").ClearContents()
